# Daily attendance processing - 2026-01-04 15:56:40
# Swap the order of "System" and the recorder's email in the
# "Recorded By" column (G) from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" for every matching row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
